$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header labels (shifted one column to the right in meaning) ---
$ws.Range("A1").Value = "Test ID"
$ws.Range("B1").Value = "Collection Date"
$ws.Range("C1").Value = "Latitude"
$ws.Range("D1").Value = "Longitude"
$ws.Range("E1").Value = "Name"
$ws.Range("F1").Value = "Area (ha)"
$ws.Range("G1").Value = "Gender"
$ws.Range("H1").Value = "Age"
$ws.Range("I1").Value = "Address"
$ws.Range("J1").Value = "Mobile No."
$ws.Range("K1").Value = "Soil pH"
$ws.Range("L1").Value = "Nitrogen"
$ws.Range("M1").Value = "Phosphorus"
$ws.Range("N1").Value = "Potassium"
$ws.Range("O1").Value = "Electrical Conductivity"
$ws.Range("P1").Value = "Temperature"
$ws.Range("Q1").Value = "Moisture"
$ws.Range("R1").Value = "Humidity"
$ws.Range("S1").Value = "Soil Health Score"
$ws.Range("T1").Value = "Recommendations"
$ws.Range("U1").Value = "Fertilizer Recommendation"

# --- Row 2: data values ---
# Text-like cells (A2, B2, J2) must stay text even though they look numeric,
# so force a text format before writing, then clear the format again so the
# cell doesn't end up carrying an extra style index.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "26"
$ws.Range("A2").ClearFormats()

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "12-04-2024"
$ws.Range("B2").ClearFormats()

$ws.Range("C2").Value = 65
$ws.Range("D2").Value = 265
$ws.Range("E2").Value = "asdjhasjkdh"
$ws.Range("F2").Value = 25
$ws.Range("G2").Value = "Female"
$ws.Range("H2").Value = 54
$ws.Range("I2").Value = "asjdhkjasd"

$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "9984564654"
$ws.Range("J2").ClearFormats()

$ws.Range("K2").Value = 7
$ws.Range("L2").Value = 250
$ws.Range("M2").Value = 120
$ws.Range("N2").Value = 263
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 26
$ws.Range("Q2").Value = 65
$ws.Range("R2").Value = 53
$ws.Range("S2").Value = 0.6128933149770953
$ws.Range("T2").Value = "Grow vegetables (tomato, brinjal, chili), fruits (mango, banana, citrus), and cash crops (sugarcane, tobacco)."
$ws.Range("U2").Value = "Apply organic amendments like compost (1-2 tonnes/ha) or vermicompost (0.5-1 tonne/ha). Follow integrated nutrient management practices. Apply chemical fertilizers like urea (0.08-0.12 tonnes/ha), single superphosphate (0.06-0.09 tonnes/ha), and muriate of potash (0.04-0.06 tonnes/ha) as per soil test recommendations and crop requirements."
